$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply row permutations (cyclic shifts) using a scratch row far below the data range.
# Destination ranges are cleared before each copy because the runtime's Range.Copy
# does not blank out a destination cell when the corresponding source cell is empty.
# Using an A:O range (instead of a full Rows() reference) keeps the worksheet's used-range
# limited to the real data columns once the scratch row is cleared again.
# Cycle: 21 -> 22 -> 23 -> 21
$ws.Range("A300:O300").ClearContents()
$ws.Range("A21:O21").Copy($ws.Range("A300:O300"))
$ws.Range("A21:O21").ClearContents()
$ws.Range("A22:O22").Copy($ws.Range("A21:O21"))
$ws.Range("A22:O22").ClearContents()
$ws.Range("A23:O23").Copy($ws.Range("A22:O22"))
$ws.Range("A23:O23").ClearContents()
$ws.Range("A300:O300").Copy($ws.Range("A23:O23"))
$ws.Range("A300:O300").ClearContents()

# Cycle: 27 -> 29 -> 27
$ws.Range("A300:O300").ClearContents()
$ws.Range("A27:O27").Copy($ws.Range("A300:O300"))
$ws.Range("A27:O27").ClearContents()
$ws.Range("A29:O29").Copy($ws.Range("A27:O27"))
$ws.Range("A29:O29").ClearContents()
$ws.Range("A300:O300").Copy($ws.Range("A29:O29"))
$ws.Range("A300:O300").ClearContents()

# Cycle: 28 -> 30 -> 28
$ws.Range("A300:O300").ClearContents()
$ws.Range("A28:O28").Copy($ws.Range("A300:O300"))
$ws.Range("A28:O28").ClearContents()
$ws.Range("A30:O30").Copy($ws.Range("A28:O28"))
$ws.Range("A30:O30").ClearContents()
$ws.Range("A300:O300").Copy($ws.Range("A30:O30"))
$ws.Range("A300:O300").ClearContents()

# Cycle: 31 -> 33 -> 32 -> 31
$ws.Range("A300:O300").ClearContents()
$ws.Range("A31:O31").Copy($ws.Range("A300:O300"))
$ws.Range("A31:O31").ClearContents()
$ws.Range("A33:O33").Copy($ws.Range("A31:O31"))
$ws.Range("A33:O33").ClearContents()
$ws.Range("A32:O32").Copy($ws.Range("A33:O33"))
$ws.Range("A32:O32").ClearContents()
$ws.Range("A300:O300").Copy($ws.Range("A32:O32"))
$ws.Range("A300:O300").ClearContents()

# Cycle: 46 -> 47 -> 46
$ws.Range("A300:O300").ClearContents()
$ws.Range("A46:O46").Copy($ws.Range("A300:O300"))
$ws.Range("A46:O46").ClearContents()
$ws.Range("A47:O47").Copy($ws.Range("A46:O46"))
$ws.Range("A47:O47").ClearContents()
$ws.Range("A300:O300").Copy($ws.Range("A47:O47"))
$ws.Range("A300:O300").ClearContents()

# Cycle: 49 -> 52 -> 51 -> 50 -> 49
$ws.Range("A300:O300").ClearContents()
$ws.Range("A49:O49").Copy($ws.Range("A300:O300"))
$ws.Range("A49:O49").ClearContents()
$ws.Range("A52:O52").Copy($ws.Range("A49:O49"))
$ws.Range("A52:O52").ClearContents()
$ws.Range("A51:O51").Copy($ws.Range("A52:O52"))
$ws.Range("A51:O51").ClearContents()
$ws.Range("A50:O50").Copy($ws.Range("A51:O51"))
$ws.Range("A50:O50").ClearContents()
$ws.Range("A300:O300").Copy($ws.Range("A50:O50"))
$ws.Range("A300:O300").ClearContents()

# Cycle: 57 -> 58 -> 59 -> 60 -> 61 -> 62 -> 57
$ws.Range("A300:O300").ClearContents()
$ws.Range("A57:O57").Copy($ws.Range("A300:O300"))
$ws.Range("A57:O57").ClearContents()
$ws.Range("A58:O58").Copy($ws.Range("A57:O57"))
$ws.Range("A58:O58").ClearContents()
$ws.Range("A59:O59").Copy($ws.Range("A58:O58"))
$ws.Range("A59:O59").ClearContents()
$ws.Range("A60:O60").Copy($ws.Range("A59:O59"))
$ws.Range("A60:O60").ClearContents()
$ws.Range("A61:O61").Copy($ws.Range("A60:O60"))
$ws.Range("A61:O61").ClearContents()
$ws.Range("A62:O62").Copy($ws.Range("A61:O61"))
$ws.Range("A62:O62").ClearContents()
$ws.Range("A300:O300").Copy($ws.Range("A62:O62"))
$ws.Range("A300:O300").ClearContents()

# Update timestamp column (O) for all data rows to reflect the new scrape time
for ($r = 2; $r -le 92; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-08-10 20:57:52"
}